$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1643.709068041747
$ws.Range("G3").Value = 2188.624045675153
$ws.Range("G4").Value = 3821.039774970653
$ws.Range("G5").Value = 2603.353473627419
$ws.Range("G6").Value = 1842.718555186179
$ws.Range("G7").Value = 3156.650941581288
$ws.Range("G8").Value = 5628.683186472136
$ws.Range("G9").Value = 8629.47258591945
$ws.Range("G10").Value = 13445.95221104066
$ws.Range("G11").Value = 14652.00799689555
$ws.Range("G12").Value = 31469.34071054289
$ws.Range("G13").Value = 18235.43065246696
$ws.Range("G14").Value = 23451.88215127289
$ws.Range("G15").Value = 29955.02158337412
$ws.Range("G16").Value = 31861.42806490845
$ws.Range("G17").Value = 24550.79320034093
$ws.Range("G18").Value = 24102.44414139766
$ws.Range("G19").Value = 8054.112808516474
$ws.Range("G20").Value = 1788.874375482555
$ws.Range("G21").Value = 10396.02590203876
$ws.Range("G22").Value = 11458.1792950775
